$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.658.55"
$ws.Range("E2").Value = "  +1.41%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.578.81"
$ws.Range("E3").Value = "  -0.43%  "

$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.51"
$ws.Range("E5").Value = "  -0.05%  "

$ws.Range("E6").Value = "  +0.05%  "

$ws.Range("E7").Value = "  +0.16%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.80"
$ws.Range("E8").Value = "  +1.86%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "24.07"
$ws.Range("E9").Value = "  +0.54%  "

$ws.Range("E10").Value = "  -1.09%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0593"
$ws.Range("E11").Value = "  -0.81%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0890"
$ws.Range("E12").Value = "  +0.40%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.802.38"
$ws.Range("E13").Value = "  -0.58%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.574.77"
$ws.Range("E14").Value = "  -0.68%  "

$ws.Range("E15").Value = "  -1.37%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "28.653.21"
$ws.Range("E16").Value = "  +1.35%  "

$ws.Range("E17").Value = "  -1.53%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "62.40"
$ws.Range("E18").Value = "  -1.21%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "231.32"
$ws.Range("E19").Value = "  +1.74%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.38"
$ws.Range("E20").Value = "  -1.02%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0693"
$ws.Range("E21").Value = "  -1.84%  "

$ws.Range("E22").Value = "  +0.15%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.90"
$ws.Range("E23").Value = "  -4.07%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.18"
$ws.Range("E24").Value = "  -1.45%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.07"
$ws.Range("E25").Value = "  +6.29%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.72"
$ws.Range("E26").Value = "  -0.22%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.04"
$ws.Range("E27").Value = "  -0.65%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.47"
$ws.Range("E28").Value = "  -1.43%  "

$ws.Range("E29").Value = "  -2.13%  "

$ws.Range("E30").Value = "  +0.12%  "

$ws.Range("E31").Value = "  +2.35%  "

$ws.Range("E32").Value = "  -2.04%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.22"
$ws.Range("E33").Value = "  -0.70%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.12"
$ws.Range("E34").Value = "  -1.27%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.399.32"
$ws.Range("E35").Value = "  +0.09%  "

$ws.Range("E36").Value = "  +3.76%  "

$ws.Range("E37").Value = "  -3.27%  "

$ws.Range("E38").Value = "  +0.84%  "

$ws.Range("E39").Value = "  +3.58%  "

$ws.Range("E40").Value = "  -0.54%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.525"
$ws.Range("E41").Value = "  -2.78%  "

$ws.Range("E42").Value = "  +0.20%  "

$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.795"
$ws.Range("E43").Value = "  -1.82%  "

$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.89"
$ws.Range("E44").Value = "  +0.81%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0464"
$ws.Range("E45").Value = "  +0.33%  "

$ws.Range("E46").Value = "  -1.81%  "

$ws.Range("E47").Value = "  -1.88%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "63.18"
$ws.Range("E48").Value = "  -1.55%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.714.43"
$ws.Range("E49").Value = "  -0.43%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "86.77"
$ws.Range("E50").Value = "  +0.02%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₆0103"
$ws.Range("E51").Value = "  -0.82%  "
